# Fix data type detection: "Jennifer" (s-003) was misfiled on the "A" sheet
# with a Score of "A", but her actual score is "A+". Move her row from the
# "A" sheet to the "A+" sheet (in the correct sorted position), updating the
# Score value to match.

$wb = $excel.ActiveWorkbook

$wsA  = $wb.Worksheets.Item("A")
$wsAp = $wb.Worksheets.Item("A+")

# Capture Jennifer's row (row 3 on sheet "A") before it shifts.
$studentId = $wsA.Cells.Item(3, 1).Value2
$name      = $wsA.Cells.Item(3, 2).Value2
$class     = $wsA.Cells.Item(3, 3).Value2

# Remove that row from sheet "A"; rows below shift up.
$wsA.Rows.Item(3).Delete()

# Insert a new row on sheet "A+" before the current row 3 (Jessica),
# shifting her down, and fill it with Jennifer's data + corrected score.
$wsAp.Rows.Item(3).Insert()
$wsAp.Cells.Item(3, 1).Value = $studentId
$wsAp.Cells.Item(3, 2).Value = $name
$wsAp.Cells.Item(3, 3).Value = $class
$wsAp.Cells.Item(3, 4).Value = "A+"
